$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 166056.6276821388
$ws.Range("F3").Value = 124542.4707616041
$ws.Range("G3").Value = 207570.7846026735
$ws.Range("E5").Value = 173974.9108291015
$ws.Range("F5").Value = 130481.1831218261
$ws.Range("G5").Value = 217468.6385363768
$ws.Range("E7").Value = 188562.5811801447
$ws.Range("F7").Value = 141421.9358851085
$ws.Range("G7").Value = 235703.2264751808
$ws.Range("E9").Value = 199938.7459048583
$ws.Range("F9").Value = 149954.0594286438
$ws.Range("G9").Value = 249923.4323810729
$ws.Range("E11").Value = 220427.0644893872
$ws.Range("F11").Value = 165320.2983670404
$ws.Range("G11").Value = 275533.830611734
$ws.Range("E13").Value = 233247.7237263739
$ws.Range("F13").Value = 174935.7927947804
$ws.Range("G13").Value = 291559.6546579674
$ws.Range("E15").Value = 244880.4889476037
$ws.Range("F15").Value = 183660.3667107028
$ws.Range("G15").Value = 306100.6111845047
$ws.Range("E16").Value = 22.77673374192634
$ws.Range("F16").Value = 17.08255030644475
$ws.Range("G16").Value = 28.47091717740792
$ws.Range("E17").Value = 260570.7137670146
$ws.Range("F17").Value = 195428.0353252609
$ws.Range("G17").Value = 325713.3922087682
$ws.Range("E18").Value = 334.8340024523636
$ws.Range("F18").Value = 251.1255018392727
$ws.Range("G18").Value = 418.5425030654545
$ws.Range("E19").Value = 273191.5798623238
$ws.Range("F19").Value = 204893.6848967428
$ws.Range("G19").Value = 341489.4748279047
$ws.Range("E20").Value = 587.0424977868289
$ws.Range("F20").Value = 440.2818733401216
$ws.Range("G20").Value = 733.8031222335361
$ws.Range("E21").Value = 286220.719092118
$ws.Range("F21").Value = 214665.5393190885
$ws.Range("G21").Value = 357775.8988651475
